# Apply cryptos list update (prices/volumes refreshed, SuiNetwork/Litecoin swapped,
# WrappedeETH inserted, VeChain removed, intervening rows shifted)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price cells so Excel does not
# auto-convert them to numbers (they must stay text, matching the source data).
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

# Set updated cell values
$ws.Cells.Item(2, 4).Value = '61.536.98'
$ws.Cells.Item(2, 5).Value = '  +1.60%  '
$ws.Cells.Item(3, 4).Value = '2.383.31'
$ws.Cells.Item(3, 5).Value = '  +1.34%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '552.42'
$ws.Cells.Item(5, 5).Value = '  +2.56%  '
$ws.Cells.Item(6, 4).Value = '140.05'
$ws.Cells.Item(6, 5).Value = '  +3.17%  '
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(8, 5).Value = '  +0.81%  '
$ws.Cells.Item(9, 4).Value = '2.384.90'
$ws.Cells.Item(9, 5).Value = '  +1.30%  '
$ws.Cells.Item(10, 5).Value = '  +4.50%  '
$ws.Cells.Item(11, 4).Value = '0.157'
$ws.Cells.Item(11, 5).Value = '  +1.85%  '
$ws.Cells.Item(12, 4).Value = '5.36'
$ws.Cells.Item(12, 5).Value = '  +2.62%  '
$ws.Cells.Item(13, 4).Value = '0.354'
$ws.Cells.Item(13, 5).Value = '  +4.31%  '
$ws.Cells.Item(14, 4).Value = '25.68'
$ws.Cells.Item(14, 5).Value = '  +4.86%  '
$ws.Cells.Item(15, 5).Value = '  +7.48%  '
$ws.Cells.Item(16, 4).Value = '2.816.59'
$ws.Cells.Item(16, 5).Value = '  +1.49%  '
$ws.Cells.Item(17, 4).Value = '61.490.23'
$ws.Cells.Item(17, 5).Value = '  +1.83%  '
$ws.Cells.Item(18, 4).Value = '2.385.46'
$ws.Cells.Item(18, 5).Value = '  +1.55%  '
$ws.Cells.Item(19, 4).Value = '10.97'
$ws.Cells.Item(19, 5).Value = '  +3.95%  '
$ws.Cells.Item(20, 4).Value = '4.18'
$ws.Cells.Item(20, 5).Value = '  +3.18%  '
$ws.Cells.Item(21, 4).Value = '321.77'
$ws.Cells.Item(21, 5).Value = '  +2.59%  '
$ws.Cells.Item(22, 4).Value = '6.71'
$ws.Cells.Item(22, 5).Value = '  +2.39%  '
$ws.Cells.Item(23, 5).Value = '  +0.02%  '
$ws.Cells.Item(24, 2).Value = 'SuiNetwork'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(24, 4).Value = '1.77'
$ws.Cells.Item(24, 5).Value = '  -5.10%  '
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(25, 4).Value = '64.27'
$ws.Cells.Item(25, 5).Value = '  +2.11%  '
$ws.Cells.Item(26, 4).Value = '8.85'
$ws.Cells.Item(26, 5).Value = '  +5.25%  '
$ws.Cells.Item(27, 4).Value = '0.998'
$ws.Cells.Item(27, 5).Value = '  -0.28%  '
$ws.Cells.Item(28, 2).Value = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).Value = '2.502.17'
$ws.Cells.Item(28, 5).Value = '  +1.49%  '
$ws.Cells.Item(29, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(29, 4).Value = '8.22'
$ws.Cells.Item(29, 5).Value = '  +4.39%  '
$ws.Cells.Item(30, 2).Value = 'Bittensor'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(30, 4).Value = '520.71'
$ws.Cells.Item(30, 5).Value = '  +4.37%  '
$ws.Cells.Item(31, 2).Value = 'PEPE'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(31, 4).Value = '0.0₃0905'
$ws.Cells.Item(31, 5).Value = '  +2.27%  '
$ws.Cells.Item(32, 2).Value = 'Fetch.AI'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(32, 4).Value = '1.39'
$ws.Cells.Item(32, 5).Value = '  +1.73%  '
$ws.Cells.Item(33, 2).Value = 'Kaspa'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(33, 4).Value = '0.149'
$ws.Cells.Item(33, 5).Value = '  +3.59%  '
$ws.Cells.Item(34, 2).Value = 'PancakeSwap'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(34, 4).Value = '1.84'
$ws.Cells.Item(34, 5).Value = '  +4.03%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = '1.51'
$ws.Cells.Item(35, 5).Value = '  +0.13%  '
$ws.Cells.Item(36, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(36, 4).Value = '1.00'
$ws.Cells.Item(36, 5).Value = '  +0.09%  '
$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(37, 4).Value = '5.57'
$ws.Cells.Item(37, 5).Value = '  +7.36%  '
$ws.Cells.Item(38, 2).Value = 'NEARProtocol'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(38, 4).Value = '4.73'
$ws.Cells.Item(38, 5).Value = '  +4.22%  '
$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39, 4).Value = '1.89'
$ws.Cells.Item(39, 5).Value = '  +7.14%  '
$ws.Cells.Item(40, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(40, 4).Value = '0.378'
$ws.Cells.Item(40, 5).Value = '  +2.22%  '
$ws.Cells.Item(41, 2).Value = 'EthereumClassic'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(41, 4).Value = '18.54'
$ws.Cells.Item(41, 5).Value = '  +1.36%  '
$ws.Cells.Item(42, 2).Value = 'Monero'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).Value = '146.20'
$ws.Cells.Item(42, 5).Value = '  +5.88%  '
$ws.Cells.Item(43, 2).Value = 'USDe'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 5).Value = '  -0.02%  '
$ws.Cells.Item(44, 2).Value = 'OKB'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(44, 4).Value = '41.41'
$ws.Cells.Item(44, 5).Value = '  +3.25%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).Value = '148.91'
$ws.Cells.Item(45, 5).Value = '  +5.73%  '
$ws.Cells.Item(46, 2).Value = 'dogwifhat'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(46, 4).Value = '2.17'
$ws.Cells.Item(46, 5).Value = '  +4.11%  '
$ws.Cells.Item(47, 2).Value = 'Filecoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(47, 4).Value = '3.62'
$ws.Cells.Item(47, 5).Value = '  +3.57%  '
$ws.Cells.Item(48, 2).Value = 'Hedera'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(48, 4).Value = '0.0525'
$ws.Cells.Item(48, 5).Value = '  +3.47%  '
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).Value = '19.84'
$ws.Cells.Item(49, 5).Value = '  +2.71%  '
$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).Value = '0.584'
$ws.Cells.Item(50, 5).Value = '  +2.91%  '
$ws.Cells.Item(51, 2).Value = 'Stellar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(51, 4).Value = '0.0906'
$ws.Cells.Item(51, 5).Value = '  +1.18%  '
